# jiko.xlsx — fix the top border row of the grid ("kaamera" frame) so the
# whole first row reads 1 instead of 0, and restore the selection to BE5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1, columns B..AM were 0 and should become 1 (A1 and AN1 are already 1).
$ws.Range("B1:AM1").Value = 1

# Put the active selection back on BE5 (was BE10).
$ws.Range("BE5").Select()
